$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Cells changing from placeholder text to a real number (need number format applied) ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'

# --- Cells changing from a real number back to placeholder text ("0" / "***.*") ---
# Force text storage via a text number format, then restore the original "General"
# look-and-feel by copying formats from a cell that already uses the placeholder style.
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null

# --- Plain numeric updates (style already correct, only the figure changed) ---
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = -4.761904761904
$ws.Range("L15").Value = 42.857142857142
$ws.Range("M15").Value = 81.818181818181
$ws.Range("N15").Value = 25
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -11.111111111111
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = -38.297872340425
$ws.Range("I16").Value = 165
$ws.Range("J16").Value = 237
$ws.Range("K16").Value = -30.379746835443
$ws.Range("L16").Value = -1.785714285714
$ws.Range("M16").Value = -7.821229050279
$ws.Range("N16").Value = -76.858345021037
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 38.461538461538
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 71
$ws.Range("H17").Value = -29.577464788732
$ws.Range("I17").Value = 312
$ws.Range("J17").Value = 349
$ws.Range("K17").Value = -10.601719197707
$ws.Range("L17").Value = 22.834645669291
$ws.Range("M17").Value = 77.272727272727
$ws.Range("N17").Value = 31.645569620253
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 86
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = -27.118644067796
$ws.Range("L18").Value = -6.521739130434
$ws.Range("M18").Value = -48.502994011976
$ws.Range("N18").Value = -92.554112554112
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 24
$ws.Range("G19").Value = 97
$ws.Range("H19").Value = -40.206185567010
$ws.Range("I19").Value = 392
$ws.Range("J19").Value = 519
$ws.Range("K19").Value = -24.470134874759
$ws.Range("L19").Value = -15.879828326180
$ws.Range("M19").Value = 57.429718875502
$ws.Range("N19").Value = -51.724137931034
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -8
$ws.Range("I20").Value = 134
$ws.Range("J20").Value = 171
$ws.Range("K20").Value = -21.637426900584
$ws.Range("L20").Value = -22.093023255814
$ws.Range("M20").Value = 5.511811023622
$ws.Range("N20").Value = -89.141004862236
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -25.862068965517
$ws.Range("F21").Value = 171
$ws.Range("G21").Value = 262
$ws.Range("H21").Value = -34.732824427480
$ws.Range("I21").Value = 1113
$ws.Range("J21").Value = 1417
$ws.Range("K21").Value = -21.453775582216
$ws.Range("L21").Value = -4.708904109589
$ws.Range("M21").Value = 22.307692307692
$ws.Range("N21").Value = -73.354081876945
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -43.859649122807
$ws.Range("M22").Value = 77.777777777777
$ws.Range("C24").Value = 45
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 118
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -7.8125
$ws.Range("I24").Value = 755
$ws.Range("J24").Value = 1204
$ws.Range("K24").Value = -37.292358803986
$ws.Range("L24").Value = -32.769367764915
$ws.Range("M24").Value = 28.183361629881
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = -45.3125
$ws.Range("I25").Value = 282
$ws.Range("J25").Value = 664
$ws.Range("K25").Value = -57.530120481927
$ws.Range("L25").Value = -51.295336787564
$ws.Range("C26").Value = 32
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 28
$ws.Range("F26").Value = 98
$ws.Range("G26").Value = 114
$ws.Range("H26").Value = -14.035087719298
$ws.Range("I26").Value = 613
$ws.Range("J26").Value = 725
$ws.Range("K26").Value = -15.448275862069
$ws.Range("L26").Value = 21.386138613861
$ws.Range("M26").Value = 21.146245059288
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = -17.142857142857
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -60
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 67
$ws.Range("J28").Value = 86
$ws.Range("K28").Value = -22.093023255814
$ws.Range("L28").Value = -19.277108433734
